$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "535.52") would be
# auto-converted from text to a numeric value by Excel when assigned directly.
# The source workbook stores these as literal text (inlineStr), so force a
# "Text" number format before writing the value, then restore the default
# "Normal" style afterwards so the cell formatting is left untouched.
$textForcedCells = @("D5", "D6", "D10", "D16", "D19", "D21", "D22", "D24", "D27", "D28", "D29", "D32", "D36", "D39", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49")
foreach ($cell in $textForcedCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.124.72"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "2.518.71"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D5").Value = "535.52"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "139.47"
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "2.523.41"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "2.965.89"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "59.092.51"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "22.93"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "2.538.34"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "10.89"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "321.18"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "62.86"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "7.75"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").Value = "6.70"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").Value = "0.0₃0765"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").Value = "160.38"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "18.48"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").Value = "36.95"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "5.25"
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("D43").Value = "282.70"
$ws.Range("E43").Value = "  -4.47%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "10.88"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").Value = "0.595"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "0.0928"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").Value = "122.96"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "18.47"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("E51").Value = "  -1.55%  "

foreach ($cell in $textForcedCells) {
    $ws.Range($cell).Style = "Normal"
}
